$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix B80: was stored as text "3", should become a real number 3 ---
$ws.Cells.Item(80, 2).Value = 3

# --- 2. Add new row 81 with annotation data ---
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = "Ruilin"
$arr[0,1] = "4"
$arr[0,2] = "good"
$arr[0,3] = "APC"
$arr[0,4] = "RES"
$arr[0,5] = "b45903b9-0b88-4991-a561-27ff9794740a"
$arr[0,6] = "S1nQvfgA-_annotated.xlsx"
$arr[0,7] = "It showed some good visualization results on controlled image generation."
$ws.Range("A81:H81").Value = $arr

# B81 ("4") must stay a text value (like the original inline string), not get
# auto-converted to a number. Writing it directly always gets coerced to a
# number, so build it as text in a scratch cell (forcing text via NumberFormat
# "@"), copy/paste-values it into B81, then fully remove the scratch cell
# (via Delete, not Clear) so it leaves no trace in the used range/dimension.
$helper = $ws.Cells.Item(200, 50)
$helper.NumberFormat = "@"
$helper.Value = "4"
$helper.Copy()
$ws.Cells.Item(81, 2).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$helper.Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
